# Add new "tipo_transaccion" and "categorias" validation lists to the
# validaciones.xlsx resource sheet (Hoja1), alongside the existing
# "tipo_cuenta" / "divisas" lists, plus one extra "tipo_cuenta" entry.
#
# Cell values are written in the same order the original author entered
# them (C1, C3, C2, D1, D3, D4, D5, D6, D7, D2, D8, D9, D10, A7, C4) so the
# shared-string table is rebuilt with the same de-duplicated ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column: tipo_transaccion (C1:C4)
$ws.Range("C1").Value = "tipo_transaccion"
$ws.Range("C3").Value = "Ingreso"
$ws.Range("C2").Value = "Gasto"

# New column: categorias (D1:D10)
$ws.Range("D1").Value = "categorias"
$ws.Range("D3").Value = "Compras"
$ws.Range("D4").Value = "Vivienda"
$ws.Range("D5").Value = "Transporte"
$ws.Range("D6").Value = "Vehículos"
$ws.Range("D7").Value = "Vida y entretenimiento"
$ws.Range("D2").Value = "Comida y bebida"
$ws.Range("D8").Value = "Comunicaciones, PC"
$ws.Range("D9").Value = "Gastos financieros"
$ws.Range("D10").Value = "Inversiones"

# Extra tipo_cuenta entry
$ws.Range("A7").Value = "Otros"

# Final tipo_transaccion entry (reuses existing shared string)
$ws.Range("C4").Value = "Transferencia"

# Extra categorias rows reusing existing shared strings
$ws.Range("D11").Value = "Ingreso"
$ws.Range("D12").Value = "Otros"

# Widen new columns to match the committed layout (bestFit-style autosize)
# and set the active selection, matching the saved workbook view.
$ws.Columns.Item(3).ColumnWidth = 14.93
$ws.Columns.Item(4).ColumnWidth = 19.76
$ws.Range("C5").Select() | Out-Null
